$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder -> speaker name
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Harold Pulcher"

# Subtitle placeholder -> role / title (two runs, second run not marked dirty)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Senior Consultant and Microsoft MVP"
